# AddToCart module / logs / readme cleanup commit:
# Rewrites the "Login" testdata sheet with the new scenario rows,
# adds two new hyperlinks, resizes a couple of columns, and moves
# the saved cell selections on both sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Register"
$ws2 = $wb.Worksheets.Item(2)   # "Login" (active sheet)

# ---------------------------------------------------------------
# Sheet "Login" (sheet2) - rebuild the scenario table
# ---------------------------------------------------------------

# Drop the old row 7 entirely (table now ends at row 6).
$ws2.Rows.Item(7).Delete()

# Row 2 - Registered user scenario
$ws2.Range("A2").Value = "Aj0007@gmail.com"
$ws2.Range("B2").Value = "Aj1234"
$ws2.Range("C2").Value = "Registered User"
$ws2.Range("D2").Value = "Login successful"

# Row 3 - UnRegistered user scenario (new hyperlink on A3)
$ws2.Range("A3").Value = "Sammed@Gmail.com"
$ws2.Range("B3").Value = 123456
$ws2.Range("C3").Value = "UnRegistered User"
$ws2.Range("D3").Value = "No customer account found"
$ws2.Range("E3").Value = "//li[normalize-space()='No customer account found']"
$ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:Sammed@Gmail.com")
$ws2.Range("A3").WrapText = $true
$ws2.Range("A3").VerticalAlignment = -4108

# Row 4 - Blank username scenario
$ws2.Range("B4").Value = "Aj1234"
$ws2.Range("D4").Value = "no customer account"

# Row 5 - Blank password scenario
$ws2.Range("A5").Value = "Ajay@gmail.com"
$ws2.Range("D5").Value = "login was unsuccessful"
$ws2.Range("E5").WrapText = $true
$ws2.Range("E5").VerticalAlignment = -4108

# Row 6 - Invalid email format scenario (new hyperlink on A6, custom display text)
$ws2.Range("B6").Value = "Aj1234"
$ws2.Range("D6").Value = "Please enter a valid email address."
$ws2.Range("E6").Value = "//span[@for='Email'"
$ws2.Range("E6").WrapText = $true
$ws2.Range("E6").VerticalAlignment = -4108
$ws2.Hyperlinks.Add($ws2.Range("A6"), "mailto:Ajay.com", "", "", "email@123.com")
$ws2.Range("A6").Value = "Ajay.com"
$ws2.Range("A6").WrapText = $true
$ws2.Range("A6").VerticalAlignment = -4108

# Widen columns D and E slightly (closest width this engine can represent).
$ws2.Columns.Item(4).ColumnWidth = 16.333333333333336
$ws2.Columns.Item(5).ColumnWidth = 18.333333333333336

# ---------------------------------------------------------------
# Selections - restore Login as active sheet/cell once Register's
# own selection has been updated.
# ---------------------------------------------------------------
$ws1.Range("C2").Select()
$ws2.Activate()
$ws2.Range("C3").Select()
